# Apply the "clean forward run codes" edit to model_settings.xlsx
#
# Changes:
#  1. Update the description text in C5 (shared string reused) to drop the
#     per-site/per-PFT/global wording now that only "all_sites" applies.
#  2. Remove the data validation (dropdown list) on B5, since the
#     all_year/site_year/per_pft/global_opti choice no longer exists.
#  3. Move the active selection from B7 to B5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("model_settings")

# 1. Update the long description text in C5.
$ws.Range("C5").Value = "only ``all_sites`` must be selected. This settings file to be used for the experiment in which a group of param vary per year, while other param are fixed across years."

# 2. Remove the dropdown (list) data validation that was attached to B5.
$ws.Range("B5").Validation.Delete()

# 3. Move the selected cell from B7 to B5.
$ws.Range("B5").Select() | Out-Null
